# "Dados_Cliente" -> "Cliente": rename the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados_Cliente")
$ws.Name = "Cliente"

# Grab the existing (horizontal, one-row-per-record) values before
# rewriting the sheet into a vertical (one-row-per-field) layout.
$empresa      = $ws.Range("A2").Value()   # "Sua Empresa"
$categoria    = $ws.Range("B2").Value()   # "Calçados"
$ticketMedio  = $ws.Range("C2").Value()   # 250
$margem       = $ws.Range("D2").Value()   # 35
$faturamento  = $ws.Range("E2").Value()   # 150000
$unidades     = $ws.Range("F2").Value()   # 600
$cac          = $ws.Range("I2").Value()   # 45
$investimento = $ws.Range("J2").Value()   # 15000

# Wipe the old header+data rows (A1:J2) so only the new A:B columns remain.
$ws.Range("A1:J2").Clear()

# Rebuild as label/value pairs going down column A (labels) / B (values).
$ws.Range("A1").Value = "Empresa"
$ws.Range("B1").Value = $empresa

$ws.Range("A2").Value = "Categoria Macro"
$ws.Range("B2").Value = $categoria

$ws.Range("A3").Value = "Ticket Médio Geral"
$ws.Range("B3").Value = $ticketMedio

$ws.Range("A4").Value = "Margem Atual"
$ws.Range("B4").Value = $margem

$ws.Range("A5").Value = "Faturamento Médio 3M"
$ws.Range("B5").Value = $faturamento

$ws.Range("A6").Value = "Unidades Médias 3M"
$ws.Range("B6").Value = $unidades

$ws.Range("A7").Value = "Range Permitido"
$ws.Range("B7").Value = 0.2

$ws.Range("A8").Value = "Ticket Customizado"
$ws.Range("B8").Value = ""
# Force the otherwise-blank B8 cell to still be materialized in the
# sheet (mirrors the source file, which keeps an explicit empty cell
# here instead of omitting it).
$ws.Range("B8").Font.Bold = $false

$ws.Range("A9").Value = "CAC"
$ws.Range("B9").Value = $cac

$ws.Range("A10").Value = "Investimento Mkt"
$ws.Range("B10").Value = $investimento
